# AIP-392 AIP-466 Updated Test Data for 2 Iteration and for group mask
#
# Expands the single data row into two rows (row2, row3) and splits the
# single "GroupMaskID" column (for each of the two devices) into four
# separate "GrpN_GroupMaskIDx" columns, shifting all following columns to
# the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1 - headers
# ---------------------------------------------------------------------
$ws.Range("A1").Value  = "Username"
$ws.Range("B1").Value  = "Password"
$ws.Range("C1").Value  = "DeviceName1"
$ws.Range("D1").Value  = "DeviceType1"
$ws.Range("E1").Value  = "DeviceIPAdd1"
$ws.Range("F1").Value  = "DeviceSerialNo1"
$ws.Range("G1").Value  = "DeviceName2"
$ws.Range("H1").Value  = "DeviceType2"
$ws.Range("I1").Value  = "DeviceIPAdd2"
$ws.Range("J1").Value  = "DeviceSerialNo2"
$ws.Range("K1").Value  = "PrefaultTime1"
$ws.Range("L1").Value  = "PostFaultTime1"
$ws.Range("M1").Value  = "MaxDFR1"
$ws.Range("N1").Value  = "UDPPortNumber1"
$ws.Range("O1").Value  = "Grp1_GroupMaskID1"
$ws.Range("P1").Value  = "Grp2_GroupMaskID1"
$ws.Range("Q1").Value  = "Grp3_GroupMaskID1"
$ws.Range("R1").Value  = "Grp4_GroupMaskID1"
$ws.Range("S1").Value  = "Compatibility1"
$ws.Range("T1").Value  = "PrefaultTime2"
$ws.Range("U1").Value  = "PostFaultTime2"
$ws.Range("V1").Value  = "MaxDFR2"
$ws.Range("W1").Value  = "UDPPortNumber2"
$ws.Range("X1").Value  = "Grp1_GroupMaskID2"
$ws.Range("Y1").Value  = "Grp2_GroupMaskID2"
$ws.Range("Z1").Value  = "Grp3_GroupMaskID2"
$ws.Range("AA1").Value = "Grp4_GroupMaskID2"
$ws.Range("AB1").Value = "Compatibility2"
$ws.Range("AC1").Value = "TimeMaster"
$ws.Range("AD1").Value = "TimeSlave"
$ws.Range("AE1").Value = "TimeMasterClock_Setting"
$ws.Range("AF1").Value = "TimeSlave_Setting_Backup_IP"
$ws.Range("AG1").Value = "TimeSlave_Setting_PPS"
$ws.Range("AH1").Value = "NoOfManualTrigger"
$ws.Range("AI1").Value = "ExpectedRecordLength"

# ---------------------------------------------------------------------
# Row 2 - iteration 1 data
# ---------------------------------------------------------------------
$ws.Range("A2").Value  = "Admin"
$ws.Range("B2").Value  = "Admin"
$ws.Range("C2").Value  = "IND_DAU_51"
$ws.Range("D2").Value  = "IDM+18"
$ws.Range("E2").Value  = "10.75.58.51"
$ws.Range("F2").Value  = "'409026540"
$ws.Range("G2").Value  = "IND_DAU_50"
$ws.Range("H2").Value  = "IDM+18"
$ws.Range("I2").Value  = "10.75.58.50"
$ws.Range("J2").Value  = "'342167760"
$ws.Range("K2").Value  = "'200"
$ws.Range("L2").Value  = "'2000"
$ws.Range("M2").Value  = "'30000"
$ws.Range("N2").Value  = "'1025"
$ws.Range("O2").Value  = "'1"
$ws.Range("P2").Value  = "'1"
$ws.Range("Q2").Value  = "'1"
$ws.Range("R2").Value  = "'1"
$ws.Range("S2").Value  = "'1"
$ws.Range("T2").Value  = "'200"
$ws.Range("U2").Value  = "'5000"
$ws.Range("V2").Value  = "'30000"
$ws.Range("W2").Value  = "'1025"
$ws.Range("X2").Value  = "'1"
$ws.Range("Y2").Value  = "'1"
$ws.Range("Z2").Value  = "'1"
$ws.Range("AA2").Value = "'1"
$ws.Range("AB2").Value = "'1"
$ws.Range("AC2").Value = "IND_DAU_51"
$ws.Range("AD2").Value = "IND_DAU_50"
$ws.Range("AE2").Value = "Internal Clock"
$ws.Range("AF2").Value = "10.75.58.51"
$ws.Range("AG2").Value = "None"
$ws.Range("AH2").Value = "'30"
$ws.Range("AI2").Value = "'30000"

# ---------------------------------------------------------------------
# Row 3 - iteration 2 data
# ---------------------------------------------------------------------
$ws.Range("A3").Value  = "Admin"
$ws.Range("B3").Value  = "Admin"
$ws.Range("C3").Value  = "IND_DAU_51"
$ws.Range("D3").Value  = "IDM+18"
$ws.Range("E3").Value  = "10.75.58.51"
$ws.Range("F3").Value  = "'409026540"
$ws.Range("G3").Value  = "IND_DAU_50"
$ws.Range("H3").Value  = "IDM+18"
$ws.Range("I3").Value  = "10.75.58.50"
$ws.Range("J3").Value  = "'342167760"
$ws.Range("K3").Value  = "'200"
$ws.Range("L3").Value  = "'5000"
$ws.Range("M3").Value  = "'30000"
$ws.Range("N3").Value  = "'1025"
$ws.Range("O3").Value  = "'1"
$ws.Range("P3").Value  = "'1"
$ws.Range("Q3").Value  = "'1"
$ws.Range("R3").Value  = "'1"
$ws.Range("S3").Value  = "'1"
$ws.Range("T3").Value  = "'200"
$ws.Range("U3").Value  = "'2000"
$ws.Range("V3").Value  = "'30000"
$ws.Range("W3").Value  = "'1025"
$ws.Range("X3").Value  = "'1"
$ws.Range("Y3").Value  = "'1"
$ws.Range("Z3").Value  = "'1"
$ws.Range("AA3").Value = "'1"
$ws.Range("AB3").Value = "'1"
$ws.Range("AC3").Value = "IND_DAU_50"
$ws.Range("AD3").Value = "IND_DAU_51"
$ws.Range("AE3").Value = "Internal Clock"
$ws.Range("AF3").Value = "10.75.58.50"
$ws.Range("AG3").Value = "None"
$ws.Range("AH3").Value = "'30"
$ws.Range("AI3").Value = "'30000"

# ---------------------------------------------------------------------
# Re-fit column widths to the new (generally longer) header/content text,
# mirroring the "best fit" auto-sizing Excel performs whenever sheet
# content changes. (ColumnWidth is expressed in characters; OOXML stores
# width = ColumnWidth + 5/6, so we back that offset out of the desired
# stored width.)
# ---------------------------------------------------------------------
function Set-BestFitWidth($col, $storedWidth) {
    $ws.Columns($col).ColumnWidth = $storedWidth - (5/6)
}

Set-BestFitWidth "A"  10
Set-BestFitWidth "B"  9.42578125
Set-BestFitWidth "C"  13.42578125
Set-BestFitWidth "D"  12.28515625
Set-BestFitWidth "E"  13.42578125
Set-BestFitWidth "F"  15.7109375
Set-BestFitWidth "G"  13.42578125
Set-BestFitWidth "H"  12.28515625
Set-BestFitWidth "I"  13.42578125
Set-BestFitWidth "J"  15.7109375
Set-BestFitWidth "K"  13.7109375
Set-BestFitWidth "L"  14.85546875
Set-BestFitWidth "N"  16.85546875
Set-BestFitWidth "O"  16.85546875
Set-BestFitWidth "P"  16.85546875
Set-BestFitWidth "Q"  16.85546875
Set-BestFitWidth "R"  14
Set-BestFitWidth "S"  14.140625
Set-BestFitWidth "T"  13.7109375
Set-BestFitWidth "U"  14.85546875
Set-BestFitWidth "W"  16.85546875
Set-BestFitWidth "X"  14
Set-BestFitWidth "Y"  14.140625
Set-BestFitWidth "Z"  12.140625
Set-BestFitWidth "AA" 12.140625
Set-BestFitWidth "AB" 28
Set-BestFitWidth "AC" 22.140625
Set-BestFitWidth "AD" 18.7109375
Set-BestFitWidth "AE" 21.7109375
Set-BestFitWidth "AF" 28

# Selection / view state left by the author after editing the sheet.
$ws.Range("AF3").Select() | Out-Null
